# Update the date line and the division problems to the new set.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-11-20 Thursday"; new = "2025-11-21 Friday"},
    @{old = "199÷8="; new = "373÷9="},
    @{old = "286÷8="; new = "948÷8="},
    @{old = "725÷3="; new = "532÷5="},
    @{old = "803÷7="; new = "748÷4="},
    @{old = "141÷3="; new = "739÷7="},
    @{old = "821÷5="; new = "568÷5="},
    @{old = "701÷6="; new = "939÷6="},
    @{old = "790÷9="; new = "474÷4="},
    @{old = "346÷9="; new = "107÷2="},
    @{old = "170÷6="; new = "781÷9="},
    @{old = "980÷6="; new = "123÷2="},
    @{old = "903÷7="; new = "560÷3="},
    @{old = "223÷6="; new = "567÷9="},
    @{old = "426÷2="; new = "964÷3="},
    @{old = "851÷7="; new = "266÷7="},
    @{old = "542÷9="; new = "887÷4="},
    @{old = "911÷6="; new = "523÷2="},
    @{old = "115÷5="; new = "216÷4="},
    @{old = "540÷2="; new = "942÷2="},
    @{old = "509÷9="; new = "947÷6="},
    @{old = "979÷6="; new = "837÷9="},
    @{old = "756÷8="; new = "924÷3="},
    @{old = "163÷7="; new = "225÷4="},
    @{old = "771÷4="; new = "419÷3="},
    @{old = "522÷4="; new = "389÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
